$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 48: the q1..q42 (and age) answers were entered as text; convert
# them to real numbers, matching every other response row in the sheet. ---
$ws.Range("B48").Value = 34

$row48Answers = @(2,3,2,3,2,2,1,2,2,3,2,2,1,2,2,3,2,2,3,2,1,2,3,2,1,2,3,2,2,3,2,1,2,3,2,2,1,2,2,2,3,2)
$col = 7
foreach ($ans in $row48Answers) {
    $ws.Cells.Item(48, $col).Value = $ans
    $col = $col + 1
}

# --- Row 49: a brand-new survey response, appended below the last row.
# Every field (including the numeric-looking answers) is stored as text;
# a leading apostrophe forces text storage for the numeric-looking values
# without touching the cell format of the plain text fields. ---
$ws.Range("A49").Value = "2025-05-16 14:59:19"
$ws.Range("B49").Value = "'23"
$ws.Range("C49").Value = "Bali, Indonesia"
$ws.Range("D49").Value = "D3"
$ws.Range("E49").Value = "male"
$ws.Range("F49").Value = "dsada"

$row49Answers = @("'4","'3","'3","'2","'3","'2","'3","'2","'2","'2","'2","'3","'2","'1","'2","'3","'2","'3","'2","'2","'2","'3","'1","'2","'3","'2","'2","'2","'3","'2","'1","'2","'2","'3","'2","'2","'1","'2","'2","'3","'2","'3")
$col = 7
foreach ($ans in $row49Answers) {
    $ws.Cells.Item(49, $col).Value = $ans
    $col = $col + 1
}
